# Auto-generated Excel COM-interop script
# Applies updated TPM values to specific cells per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3197726666666666
$ws.Range("H2").Value = 0.959318
$ws.Range("I2").Value = 0.1074590987069417
$ws.Range("J2").Value = 0.1074590987069417
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 24.65112244836222
$ws.Range("R2").Value = 221.86010203526
$ws.Range("S2").Value = 0.0258311973215403
$ws.Range("T2").Value = 0.0258311973215403
$ws.Range("G3").Value = 0.3197726666666666
$ws.Range("H3").Value = 0.959318
$ws.Range("I3").Value = 0.1074590987069417
$ws.Range("J3").Value = 0.1074590987069417
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 32.48251941817956
$ws.Range("R3").Value = 292.342674763616
$ws.Range("S3").Value = 0.0340374914103559
$ws.Range("T3").Value = 0.03403749141035589
$ws.Range("G4").Value = 0.3197726666666666
$ws.Range("H4").Value = 0.959318
$ws.Range("I4").Value = 0.1074590987069417
$ws.Range("J4").Value = 0.1074590987069417
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 45.41628516322489
$ws.Range("R4").Value = 408.746566469024
$ws.Range("S4").Value = 0.04759040997504549
$ws.Range("T4").Value = 0.04759040997504548
$ws.Range("I5").Value = 0.7259488187057992
$ws.Range("J5").Value = 0.7259488187057991
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 166.5326941738489
$ws.Range("R5").Value = 1498.79424756464
$ws.Range("S5").Value = 0.1745047874677292
$ws.Range("T5").Value = 0.1745047874677291
$ws.Range("I6").Value = 0.7259488187057992
$ws.Range("J6").Value = 0.7259488187057991
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.2299430851215623
$ws.Range("T6").Value = 0.2299430851215622
$ws.Range("I7").Value = 0.7259488187057992
$ws.Range("J7").Value = 0.7259488187057991
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.3215009461165078
$ws.Range("T7").Value = 0.3215009461165078
$ws.Range("I8").Value = 0.1665920825872592
$ws.Range("J8").Value = 0.1665920825872591
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 38.21623181472778
$ws.Range("R8").Value = 343.94608633255
$ws.Range("S8").Value = 0.04004568258341298
$ws.Range("T8").Value = 0.04004568258341296
$ws.Range("I9").Value = 0.1665920825872592
$ws.Range("J9").Value = 0.1665920825872591
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.05276776604614159
$ws.Range("T9").Value = 0.05276776604614156
$ws.Range("I10").Value = 0.1665920825872592
$ws.Range("J10").Value = 0.1665920825872591
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.07377863395770463
$ws.Range("T10").Value = 0.07377863395770461
